# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-data refresh described in the commit diff.
# For each sheet, a set of cells in columns H-N (currentAveragePrice*, LevePrice*, LeveProfit*)
# is updated to reflect freshly fetched market data. A few cells are cleared entirely
# (no HQ/NQ profit applicable) and one cell is newly populated.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$updates = @{
    "H10" = 0
    "I10" = 0
    "K10" = 0
    "H40" = 4999
    "H58" = 1444.75
    "I58" = 417.125
    "K58" = 1251.375
    "M58" = -1101.375
    "H61" = 850.3333
    "I61" = 750.375
    "K61" = 2251.125
    "M61" = -2079.125
    "H62" = 24269.6
    "J62" = 32273.9
    "L62" = 32273.9
    "N62" = -33521.9
    "H65" = 24269.6
    "J65" = 32273.9
    "L65" = 161369.5
    "N65" = -167609.5
    "H82" = 292.5
    "I82" = 292.5
    "K82" = 877.5
    "M82" = -471.5
    "H85" = 292.5
    "I85" = 292.5
    "K85" = 877.5
    "M85" = 526.5
    "H101" = 12988777
    "I101" = 17859076
    "J101" = 1314.3334
    "K101" = 53577228
    "L101" = 3943.0002
    "M101" = -53575606
    "N101" = -7187.0002
    "H104" = 1439.6
    "I104" = 1439.6
    "K104" = 4318.799999999999
    "M104" = -2571.799999999999
    "H115" = 3650.5715
    "I115" = 431
    "K115" = 1293
    "M115" = 274
    "H129" = 1628
    "I129" = 1628
    "K129" = 4884
    "M129" = 116
    "H138" = 3656.8245
    "I138" = 1132.7222
    "K138" = 3398.1666
    "M138" = 1741.8334
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
foreach ($cellRef in @("M10")) {
    $ws.Range($cellRef).ClearContents()
}

$ws = $wb.Worksheets.Item("ARM")

$updates = @{
    "H19" = 8450
    "J19" = 0
    "L19" = 0
    "H36" = 7000
    "I36" = 7000
    "K36" = 7000
    "M36" = -6654
    "H132" = 3165.641
    "I132" = 1850.0358
    "J132" = 6514.4546
    "K132" = 5550.107400000001
    "L132" = 19543.3638
    "M132" = -3020.107400000001
    "N132" = -24603.3638
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
foreach ($cellRef in @("N19")) {
    $ws.Range($cellRef).ClearContents()
}

$ws = $wb.Worksheets.Item("BSM")

$updates = @{
    "H99" = 15914.27
    "I99" = 17240.262
    "J99" = 5748.3335
    "K99" = 17240.262
    "L99" = 5748.3335
    "M99" = -15742.262
    "N99" = -8744.333500000001
    "H105" = 3220.3333
    "J105" = 5666.6665
    "L105" = 5666.6665
    "N105" = -9160.666499999999
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("CRP")

$updates = @{
    "H35" = 3349.8
    "I35" = 3999.6667
    "K35" = 3999.6667
    "M35" = -3705.6667
    "H58" = 3803.88
    "I58" = 2595.2307
    "J58" = 5113.25
    "K58" = 2595.2307
    "L58" = 5113.25
    "M58" = -2392.2307
    "N58" = -5519.25
    "H136" = 3803.88
    "I136" = 2595.2307
    "J136" = 5113.25
    "K136" = 7785.6921
    "L136" = 15339.75
    "M136" = -5235.6921
    "N136" = -20439.75
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("CUL")

$updates = @{
    "H4" = 49690780
    "I4" = 36897464
    "J4" = 420696960
    "K4" = 110692392
    "L4" = 1262090880
    "M4" = -110692280
    "N4" = -1262091104
    "H17" = 277.64517
    "J17" = 1053.5
    "L17" = 3160.5
    "N17" = -3498.5
    "H39" = 2092.5
    "I39" = 271
    "J39" = 3914
    "K39" = 813
    "L39" = 11742
    "M39" = -519
    "N39" = -12330
    "H55" = 6281.579
    "J55" = 7247.625
    "L55" = 21742.875
    "N55" = -22096.875
    "H56" = 5582.3335
    "I56" = 5582.3335
    "K56" = 5582.3335
    "M56" = -5052.3335
    "H58" = 9375.75
    "I58" = 2501
    "K58" = 7503
    "M58" = -7375
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("GSM")

$updates = @{
    "H62" = 49000
    "J62" = 49000
    "L62" = 49000
    "N62" = -50372
    "H65" = 49000
    "J65" = 49000
    "L65" = 147000
    "N65" = -153864
    "H80" = 18331
    "I80" = 25499
    "J80" = 3995
    "K80" = 25499
    "L80" = 3995
    "M80" = -24501
    "N80" = -5991
    "H83" = 18331
    "I83" = 25499
    "J83" = 3995
    "K83" = 127495
    "L83" = 19975
    "M83" = -122503
    "N83" = -29959
    "H122" = 24850.25
    "J122" = 23601.2
    "L122" = 70803.60000000001
    "N122" = -75703.60000000001
    "H132" = 3159.0417
    "I132" = 3199.7827
    "J132" = 2222
    "K132" = 9599.348100000001
    "L132" = 6666
    "M132" = -7069.348100000001
    "N132" = -11726
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("LTW")

$updates = @{
    "H22" = 2189.818
    "I22" = 2298.4
    "J22" = 1957.1428
    "K22" = 2298.4
    "L22" = 1957.1428
    "M22" = -2003.4
    "N22" = -2547.1428
    "H27" = 2189.818
    "I27" = 2298.4
    "J27" = 1957.1428
    "K27" = 2298.4
    "L27" = 1957.1428
    "M27" = -2191.4
    "N27" = -2171.1428
    "H122" = 3760.359
    "I122" = 2702.3076
    "K122" = 8106.9228
    "M122" = -5656.9228
    "H136" = 6329.7
    "J136" = 7773.467
    "L136" = 23320.401
    "N136" = -28420.401
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$ws = $wb.Worksheets.Item("WVR")

$updates = @{
    "H14" = 0
    "I14" = 0
    "K14" = 0
    "H82" = 44066.668
    "J82" = 43600
    "L82" = 43600
    "N82" = -44366
    "H85" = 44066.668
    "J85" = 43600
    "L85" = 43600
    "N85" = -46252
    "H92" = 275025000
    "J92" = 275025000
    "L92" = 275025000
    "N92" = -275029992
    "H100" = 27217.875
    "I100" = 18655.389
    "J100" = 52905.332
    "K100" = 37310.778
    "L100" = 105810.664
    "M100" = -36769.778
    "N100" = -106892.664
    "H122" = 8377.9
    "I122" = 3957.375
    "K122" = 11872.125
    "M122" = -9422.125
    "H136" = 3278.8667
    "I136" = 2828
    "K136" = 8484
    "M136" = -5934
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
foreach ($cellRef in @("M14")) {
    $ws.Range($cellRef).ClearContents()
}
